# Insert a new row at position 53, shifting all existing rows (53-131) down
# by one. This pushes the previous last row (131) out to the new row 132
# and makes room for a brand-new weekly price record at row 53.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("53:53").Insert()

$ws.Cells.Item(53, 1).Value = 1
$ws.Cells.Item(53, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(53, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(53, 4).Value = 45012
$ws.Cells.Item(53, 5).Value = 15
$ws.Cells.Item(53, 6).Value = 100112021
$ws.Cells.Item(53, 7).Value = "Ají"
$ws.Cells.Item(53, 8).Value = "Inferno"
$ws.Cells.Item(53, 9).Value = "Primera"
$ws.Cells.Item(53, 10).Value = 150
$ws.Cells.Item(53, 11).Value = 23000
$ws.Cells.Item(53, 12).Value = 24000
$ws.Cells.Item(53, 13).Value = 23500
$ws.Cells.Item(53, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(53, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(53, 16).Value = 1567
$ws.Cells.Item(53, 17).Value = 15
$ws.Cells.Item(53, 18).Value = "Hortaliza"
